$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '71.454.68'
$ws.Cells.Item(2, 5).Value = '  -0.73%  '
$ws.Cells.Item(3, 4).Value = '2.609.57'
$ws.Cells.Item(3, 5).Value = '  -3.25%  '
$ws.Cells.Item(4, 5).Value = '  +0.22%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '584.14'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -2.15%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '172.37'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -2.13%  '
$ws.Cells.Item(7, 5).Value = '  +0.16%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.516'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -1.33%  '
$ws.Cells.Item(9, 4).Value = '2.610.94'
$ws.Cells.Item(9, 5).Value = '  -3.15%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.167'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -1.95%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.169'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +0.57%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.356'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +0.74%  '
$ws.Cells.Item(13, 5).Value = '  -2.48%  '
$ws.Cells.Item(14, 4).Value = '3.089.47'
$ws.Cells.Item(14, 5).Value = '  -2.86%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.0000187'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +0.96%  '
$ws.Cells.Item(16, 4).Value = '71.385.81'
$ws.Cells.Item(16, 5).Value = '  -0.56%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '25.53'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -2.73%  '
$ws.Cells.Item(18, 4).Value = '2.634.81'
$ws.Cells.Item(18, 5).Value = '  -2.06%  '
$ws.Cells.Item(19, 5).Value = '  -0.60%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '7.81'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -2.75%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '370.45'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +0.21%  '
$ws.Cells.Item(22, 5).Value = '  -2.48%  '
$ws.Cells.Item(23, 5).Value = '  -0.32%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '71.37'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.83%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.999'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -0.10%  '
$ws.Cells.Item(26, 5).Value = '  -3.37%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '9.33'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -5.46%  '
$ws.Cells.Item(28, 4).Value = '2.763.28'
$ws.Cells.Item(28, 5).Value = '  -2.34%  '
$ws.Cells.Item(29, 5).Value = '  +0.33%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0941'
$ws.Cells.Item(30, 5).Value = '  -2.22%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '7.92'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -1.84%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '484.64'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -4.56%  '
$ws.Cells.Item(33, 5).Value = '  +1.25%  '
$ws.Cells.Item(34, 5).Value = '  -1.38%  '
$ws.Cells.Item(35, 5).Value = '  +0.13%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '159.43'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -2.02%  '
$ws.Cells.Item(37, 5).Value = '  +5.89%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '19.19'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -1.47%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '18.90'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -0.84%  '
$ws.Cells.Item(40, 5).Value = '  -2.27%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '1.70'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -5.58%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '2.54'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -1.39%  '
$ws.Cells.Item(44, 5).Value = '  -4.05%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.324'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -2.65%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '38.85'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -0.95%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '148.77'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -4.30%  '
$ws.Cells.Item(48, 5).Value = '  -2.41%  '
$ws.Cells.Item(49, 5).Value = '  -2.28%  '
$ws.Cells.Item(50, 5).Value = '  -4.86%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.601'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -1.04%  '
